# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the "展览", "演出" and "全部类型" sheets, as generated by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12918
$ws1.Range("F4").Value = 313
$ws1.Range("F5").Value = 631
$ws1.Range("F6").Value = 206
$ws1.Range("F7").Value = 400
$ws1.Range("F8").Value = 1216

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28
$ws2.Range("G2").Value = 188
$ws2.Range("F3").Value = 11

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 28
$ws4.Range("G3").Value = 188
$ws4.Range("F4").Value = 12918
$ws4.Range("F5").Value = 313
$ws4.Range("F6").Value = 631
$ws4.Range("F7").Value = 206
$ws4.Range("F8").Value = 11
$ws4.Range("F10").Value = 400
$ws4.Range("F11").Value = 1216
